$d = $word.ActiveDocument

# 1) Replace the paragraph text "Presentation has been recorded. The file is in the Github located at:"
#    with "Presentation has been recorded. The file is located at:"
#    Use wildcard find/replace to span across the run boundaries (Github / proofErr marks).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Presentation has been recorded. The file is in the Github located at:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Presentation has been recorded. The file is located at:",
    2
) | Out-Null

# 2) Update the hyperlink display text and address.
$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.TextToDisplay = "https://drive.google.com/file/d/1MtpRa7NXa7klpLGPTbV1qwIGMf6XOdmp/view?usp=sharing"
$hyperlink.Address = "https://drive.google.com/file/d/1MtpRa7NXa7klpLGPTbV1qwIGMf6XOdmp/view?usp=sharing"
